$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) values
$ws.Range("A1").Value = "datnum"
$ws.Range("B1").Value = "datname"
$ws.Range("C1").Value = "time"
$ws.Range("D1").Value = "picklepath"
$ws.Range("E1").Value = "x_label"
$ws.Range("F1").Value = "y_label"

# Data row (row 2) values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "base"
$ws.Range("C2").Value = 1577779312.313096
$ws.Range("D2").Value = "pathtopickle"
$ws.Range("E2").Value = "xlabel"
$ws.Range("F2").Value = "ylabel"

# Copy the existing header style (from B1) onto the newly styled cells:
# A1, E1, F1 (row 1) and A2, B2 (row 2)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
